$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menus_and_forms")
$ws.Activate()

# Rename the icon/audio filepath columns to the new image/audio column names,
# removing backwards compatibility for multi-sheet bulk app translations.
$ws.Range("E1").Value = "image_en"
$ws.Range("F1").Value = "audio_en"
$ws.Range("G1").Value = "image_fra"
$ws.Range("H1").Value = "audio_fra"

# Update the active selection to match the new layout.
$ws.Range("E2").Select()
